$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 12:11 PM"

# --- "distance from Dma50" sheet: update stock order / values ---
$ws = $wb.Worksheets.Item("distance from Dma50")

$ws.Range("C2").Value = 9.5989
$ws.Range("C3").Value = 7.3337
$ws.Range("C4").Value = 6.5539
$ws.Range("C5").Value = 5.3422
$ws.Range("C6").Value = 5.2626
$ws.Range("C7").Value = 5.064
$ws.Range("C8").Value = 4.4135
$ws.Range("C9").Value = 4.3862
$ws.Range("C10").Value = 3.9007
$ws.Range("C11").Value = 3.7283

# Rows 12 and 13 swap stock names, with new values
$ws.Range("B12").Value = "CNXMIDCAP"
$ws.Range("C12").Value = 3.4128
$ws.Range("B13").Value = "NIFTYFINSERVICE"
$ws.Range("C13").Value = 3.4041

$ws.Range("C14").Value = 3.0813
$ws.Range("C15").Value = 3.0573
$ws.Range("C16").Value = 2.9726
$ws.Range("C17").Value = 2.8821
$ws.Range("C18").Value = 2.8418
$ws.Range("C19").Value = 2.7851
$ws.Range("C20").Value = 2.4001
$ws.Range("C21").Value = 2.3225
$ws.Range("C22").Value = 1.4105
$ws.Range("C23").Value = 1.3283
$ws.Range("C24").Value = 1.3242
$ws.Range("C25").Value = 1.0993
$ws.Range("C26").Value = 0.9802999999999999
$ws.Range("C27").Value = 0.8948
$ws.Range("C28").Value = 0.5919
$ws.Range("C29").Value = 0.4064
$ws.Range("C30").Value = -2.1175
